# Apply league data refresh for "Poland Ekstraklasa" sheet
# (odds/results update + row reordering for 2024-02-16 20:23 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 147
$ws.Range("B147").Value2 = 5460884
$ws.Range("F147").Value2 = "Miedz Legnica"
$ws.Range("G147").Value2 = "Gornik Zabrze"
$ws.Range("H147").Value2 = 0
$ws.Range("J147").Value2 = "D"
$ws.Range("K147").Value2 = 3.6
$ws.Range("L147").Value2 = 3.5
$ws.Range("M147").Value2 = 1.909
$ws.Range("N147").Value2 = 3
$ws.Range("O147").Value2 = 3.5
$ws.Range("P147").Value2 = 2.1
$ws.Range("Q147").Value2 = 0.25
$ws.Range("R147").Value2 = 1.95
$ws.Range("S147").Value2 = 1.9
$ws.Range("T147").Value2 = 2.75
$ws.Range("W147").Value2 = -1
$ws.Range("X147").Value2 = 2.5
$ws.Range("Z147").Value2 = 0.475
$ws.Range("AA147").Value2 = -0.5

# Row 149
$ws.Range("B149").Value2 = 5456594
$ws.Range("F149").Value2 = "Rakow Czestochowa"
$ws.Range("G149").Value2 = "Zaglebie Lubin"
$ws.Range("H149").Value2 = 1
$ws.Range("I149").Value2 = 1
$ws.Range("J149").Value2 = "D"
$ws.Range("K149").Value2 = 1.444
$ws.Range("L149").Value2 = 4.5
$ws.Range("M149").Value2 = 5.75
$ws.Range("N149").Value2 = 1.3
$ws.Range("O149").Value2 = 5.25
$ws.Range("P149").Value2 = 7
$ws.Range("Q149").Value2 = -1.5
$ws.Range("R149").Value2 = 1.9
$ws.Range("S149").Value2 = 1.95
$ws.Range("T149").Value2 = 3
$ws.Range("U149").Value2 = 1.9
$ws.Range("V149").Value2 = 1.95
$ws.Range("W149").Value2 = -1
$ws.Range("X149").Value2 = 4.25
$ws.Range("Z149").Value2 = -1
$ws.Range("AA149").Value2 = 0.95
$ws.Range("AB149").Value2 = -1
$ws.Range("AC149").Value2 = 0.95

# Row 150
$ws.Range("B150").Value2 = 5456603
$ws.Range("F150").Value2 = "Lech Poznan"
$ws.Range("G150").Value2 = "Jagiellonia Bialystok"
$ws.Range("H150").Value2 = 2
$ws.Range("I150").Value2 = 0
$ws.Range("K150").Value2 = 1.363
$ws.Range("L150").Value2 = 4.75
$ws.Range("M150").Value2 = 6.5
$ws.Range("N150").Value2 = 1.222
$ws.Range("O150").Value2 = 5.5
$ws.Range("P150").Value2 = 8
$ws.Range("Q150").Value2 = -1.75
$ws.Range("R150").Value2 = 1.925
$ws.Range("S150").Value2 = 1.925
$ws.Range("T150").Value2 = 3.25
$ws.Range("U150").Value2 = 1.95
$ws.Range("V150").Value2 = 1.9
$ws.Range("W150").Value2 = 0.222
$ws.Range("Z150").Value2 = 0.4625
$ws.Range("AA150").Value2 = -0.5
$ws.Range("AB150").Value2 = -1
$ws.Range("AC150").Value2 = 0.8999999999999999

# Row 151
$ws.Range("B151").Value2 = 5467427
$ws.Range("F151").Value2 = "Stal Mielec"
$ws.Range("G151").Value2 = "Warta Poznan"
$ws.Range("I151").Value2 = 0
$ws.Range("J151").Value2 = "H"
$ws.Range("K151").Value2 = 2.375
$ws.Range("L151").Value2 = 3.2
$ws.Range("M151").Value2 = 2.8
$ws.Range("N151").Value2 = 2.6
$ws.Range("O151").Value2 = 3.1
$ws.Range("P151").Value2 = 2.625
$ws.Range("Q151").Value2 = 0
$ws.Range("R151").Value2 = 1.925
$ws.Range("S151").Value2 = 1.925
$ws.Range("T151").Value2 = 2.25
$ws.Range("U151").Value2 = 1.975
$ws.Range("V151").Value2 = 1.875
$ws.Range("W151").Value2 = 1.6
$ws.Range("X151").Value2 = -1
$ws.Range("Z151").Value2 = 0.925
$ws.Range("AA151").Value2 = -1
$ws.Range("AC151").Value2 = 0.875

# Row 152
$ws.Range("B152").Value2 = 5461474
$ws.Range("F152").Value2 = "Legia Warsaw"
$ws.Range("G152").Value2 = "Slask Wroclaw"
$ws.Range("H152").Value2 = 3
$ws.Range("I152").Value2 = 1
$ws.Range("K152").Value2 = 1.7
$ws.Range("L152").Value2 = 3.8
$ws.Range("M152").Value2 = 4
$ws.Range("N152").Value2 = 1.833
$ws.Range("O152").Value2 = 3.8
$ws.Range("P152").Value2 = 3.4
$ws.Range("Q152").Value2 = -0.5
$ws.Range("R152").Value2 = 1.825
$ws.Range("S152").Value2 = 2.025
$ws.Range("T152").Value2 = 2.75
$ws.Range("U152").Value2 = 1.9
$ws.Range("V152").Value2 = 1.95
$ws.Range("W152").Value2 = 0.833
$ws.Range("Z152").Value2 = 0.825
$ws.Range("AA152").Value2 = -1
$ws.Range("AB152").Value2 = 0.8999999999999999
$ws.Range("AC152").Value2 = -1

# Row 154
$ws.Range("B154").Value2 = 5465446
$ws.Range("F154").Value2 = "Cracovia Krakow"
$ws.Range("G154").Value2 = "Wisla Plock"
$ws.Range("H154").Value2 = 3
$ws.Range("J154").Value2 = "H"
$ws.Range("K154").Value2 = 2.15
$ws.Range("M154").Value2 = 2.875
$ws.Range("N154").Value2 = 2.25
$ws.Range("O154").Value2 = 3.6
$ws.Range("P154").Value2 = 2.7
$ws.Range("Q154").Value2 = -0.25
$ws.Range("R154").Value2 = 2.05
$ws.Range("S154").Value2 = 1.75
$ws.Range("T154").Value2 = 2.5
$ws.Range("U154").Value2 = 1.825
$ws.Range("V154").Value2 = 2.025
$ws.Range("W154").Value2 = 1.25
$ws.Range("X154").Value2 = -1
$ws.Range("Z154").Value2 = 1.05
$ws.Range("AA154").Value2 = -1
$ws.Range("AB154").Value2 = 0.825
$ws.Range("AC154").Value2 = -1

# Row 332
$ws.Range("H332").Value2 = 0
$ws.Range("I332").Value2 = 4
$ws.Range("J332").Value2 = "A"
$ws.Range("N332").Value2 = 3.4
$ws.Range("O332").Value2 = 3.25
$ws.Range("P332").Value2 = 2.15
$ws.Range("R332").Value2 = 2
$ws.Range("S332").Value2 = 1.85
$ws.Range("W332").Value2 = -1
$ws.Range("X332").Value2 = -1
$ws.Range("Y332").Value2 = 1.15
$ws.Range("Z332").Value2 = -1
$ws.Range("AA332").Value2 = 0.8500000000000001
$ws.Range("AB332").Value2 = 0.925
$ws.Range("AC332").Value2 = -1

# Row 333
$ws.Range("H333").Value2 = 0
$ws.Range("I333").Value2 = 1
$ws.Range("J333").Value2 = "A"
$ws.Range("N333").Value2 = 1.615
$ws.Range("P333").Value2 = 5.75
$ws.Range("R333").Value2 = 1.85
$ws.Range("S333").Value2 = 2
$ws.Range("W333").Value2 = -1
$ws.Range("X333").Value2 = -1
$ws.Range("Y333").Value2 = 4.75
$ws.Range("Z333").Value2 = -1
$ws.Range("AA333").Value2 = 1
$ws.Range("AB333").Value2 = -1
$ws.Range("AC333").Value2 = 0.925

# Row 334
$ws.Range("N334").Value2 = 2.625
$ws.Range("P334").Value2 = 3
$ws.Range("R334").Value2 = 1.775
$ws.Range("S334").Value2 = 2.1

# Row 335
$ws.Range("R335").Value2 = 1.8
$ws.Range("S335").Value2 = 2.05

# Row 336
$ws.Range("U336").Value2 = 1.9
$ws.Range("V336").Value2 = 1.95

# Row 339
$ws.Range("N339").Value2 = 3.25
$ws.Range("P339").Value2 = 2.25
$ws.Range("R339").Value2 = 1.875
$ws.Range("S339").Value2 = 1.975
$ws.Range("U339").Value2 = 2.05
$ws.Range("V339").Value2 = 1.8

